{"js": "// Remove the trailing \"empty paragraph\" + \"Ver no Jupiter...\" paragraph +\n// \"\u00a9 2020 ... Attribution\" paragraph that followed the last bibliography\n// entry (\"...Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais ... ATLAS\").\n// The paragraph right after those three (another empty \"Normal\" paragraph,\n// immediately preceding the page-break paragraph) is left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the anchor paragraph (\"...ATLAS\") that precedes the block we must\n// delete, then remove: the empty paragraph right after it, plus the two\n// text paragraphs that follow.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t && t.indexOf(\"Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais\") !== -1 &&\n      t.indexOf(\"Bertero\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate the 'Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais' paragraph.\");\n}\n\nconst toDelete = [];\n// The empty paragraph immediately after the anchor.\nif (items[anchorIndex + 1] && items[anchorIndex + 1].text === \"\") {\n  toDelete.push(items[anchorIndex + 1]);\n}\n// The two text paragraphs that must be removed, matched by their content.\nfor (let i = anchorIndex + 2; i < items.length; i++) {\n  const t = items[i].text;\n  if (targetTexts.indexOf(t) !== -1) {\n    toDelete.push(items[i]);\n  }\n  if (toDelete.length >= 3) break;\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n# the \"\u00a9 2020 . Contact: ... Creative Commons Attribution\" paragraph that\n# follows it, and the blank paragraph that separates them from the last\n# bibliography entry (\"...Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais\n# ... ATLAS\"). The blank paragraph right before the page-break paragraph\n# (further down) is left untouched.\n\n$d = $word.ActiveDocument\n\n$markers = @()\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    $t = $p.Range.Text\n    if ($t -like \"*Ver no Jupiter Salvar em pdf Salvar em docx*\") {\n        $markers = $markers + $i\n    }\n    if ($t -like \"*Contact: luizeleno@usp.br*\") {\n        $markers = $markers + $i\n    }\n}\n\nif ($markers.Count -eq 0) {\n    Write-Output \"Target paragraphs not found; nothing deleted.\"\n} else {\n    $first = $markers[0]\n\n    # Include the blank paragraph immediately preceding the block, if it is\n    # indeed empty (just the paragraph mark, length 1).\n    $precedingIndex = $first - 1\n    if ($precedingIndex -ge 1) {\n        $precedingText = $d.Paragraphs.Item($precedingIndex).Range.Text\n        if ($precedingText.Length -eq 1) {\n            $markers = @($precedingIndex) + $markers\n        }\n    }\n\n    # Delete from the highest index down to the lowest so earlier indices\n    # stay valid while we work.\n    $sorted = $markers | Sort-Object -Descending\n    foreach ($idx in $sorted) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
